$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76 currently carries the "last row" date-only format (A76).
# Since a new last row (77) is being appended, row 76's date format
# reverts to the standard datetime format used by all other data rows.
$ws.Range("A76").NumberFormat = $ws.Range("A75").NumberFormat

# Add the new row 77 with the daily update values.
$ws.Range("A77").Value = 45664
$ws.Range("A77").NumberFormat = "YYYY-MM-DD"

$ws.Range("B77").Value = 180
$ws.Range("C77").Value = 178
$ws.Range("D77").Value = 180
